$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The template has a "blank separator" row right after the last filled diary
# entry (row 39, using styles 27/28/29) followed by plain blank rows (styles
# 30/31/32). We are filling in 4 new diary entries in rows 39-42, so the
# blank-separator formatting needs to move down to row 43.

# Step 1: preserve the blank-separator formatting (currently row 39) onto row 43
$ws.Range("A39:G39").Copy()
$ws.Range("A43:G43").PasteSpecial(-4122)  # xlPasteFormats

# Step 2: copy the filled-entry formatting (currently row 38) onto rows 39:42
$ws.Range("A38:G38").Copy()
$ws.Range("A39:G42").PasteSpecial(-4122)  # xlPasteFormats

# Step 3: set the row heights for the new entries (and restore row 43's
# default blank height)
$ws.Rows.Item(39).RowHeight = 199
$ws.Rows.Item(40).RowHeight = 157
$ws.Rows.Item(41).RowHeight = 151.2
$ws.Rows.Item(42).RowHeight = 129
$ws.Rows.Item(43).RowHeight = 15.5

# Step 4: Row 39 - Lecture 8 (Thu 2/27/2020)
$ws.Range("A39").Value = 43888
$ws.Range("B39").Value = "17:00 - 19:50"
$ws.Range("C39").Value = "None"
$ws.Range("D39").Value = "Looking forward to discuss our assignment on contribution standards and probably how to make the pull request. Don’t know what is going to be covered this lecture! Looking forward to the guest as usual. "
$ws.Range("E39").Value = "We discussed about the architecture and it was great to know that everyone had problems in finding the architecture or settling on a pure style. We faced a lot of difficulties in settling on the report. We learnt more about KEP. The one on invest and save time later is great.Especially because we are all changing majors I guess it is important to understand the root of all the courses than to rush. We learnt about software design patterns and the ducks example helped us in understanding the underlying concept. We had an amazing talk by Alberto."
$ws.Range("F39").Value = "Initial explanation about DBH and conference helped in understanding design pattern. We were not able to do the practical part so I am guessing the homework will be harder. Finally, listening to Alberto about Astrophysics was great because I could related more as I have done research in Computational Chemistry before with Simulation software that used Fortran. It was great to listen to some Physics again! He is very passionate about his work and it is great to know the various fields that we can work in or contribute to. "
$ws.Range("G39").Value = "Feeling tired!"

# Step 5: Row 40 - Team meeting (Sun 3/1/2020)
$ws.Range("A40").Value = 43891
$ws.Range("B40").Value = "18:00-20:00"
$ws.Range("C40").Value = "Team"
$ws.Range("D40").Value = "Finalize the issue to solve, Identity two patterns"
$ws.Range("E40").Value = "Yitian suggested we address the issue about the titles as it is relatively new and that it was tagged as good for first time contributors. We agreed on this issue after skimming through the others. We came up with basic strategy of the problem and we plan to implement this modification for the pull request. Upon looking for design patterns, we managed to look for keywords like adapter, builder, etc."
$ws.Range("F40").Value = "I realized that contributing to a project is not that bad after all and we don’t have to know each and every part because we looked into the issue and it said exactly what the problem was. We used the knowledge from previous classes to look for beacons and infant managed to find the algorithm that caused this error and to our surprise this was a concise file. "
$ws.Range("G40").Value = "Feeling smart!"

# Step 6: Row 41 - Team meeting (Mon 3/2/2020)
$ws.Range("A41").Value = 43892
$ws.Range("B41").Value = "21:00 - 23:00"
$ws.Range("C41").Value = "Team"
$ws.Range("D41").Value = "Submit pull request, Identify three patterns"
$ws.Range("E41").Value = "We agreed on the pull request submission and made it as per standards we had found in the previous assignment. We found enum iterator pattern, and it was relatively harder to find the other prototype and factory method"
$ws.Range("F41").Value = "We are satisfied that we were able to contribute to the open source community and sure that there will be more comments on the issue. Hope to hear back them in a day as this is a highly active project. While submitting this request, we also came up with a solution for the next week’s assignment. The design patterns were interesting to read and analyse. "
$ws.Range("G41").Value = "Feeling motivated to keep track of JabRef and identify minor issues to contribute!"

# Step 7: Row 42 - Team meeting (Wed 3/4/2020)
$ws.Range("A42").Value = 43894
$ws.Range("B42").Value = "21:00-23:30"
$ws.Range("C42").Value = "Team"
$ws.Range("D42").Value = " Settle on the design patterns report"
$ws.Range("E42").Value = "We finished writing the complete report for the assignment. While we checked the pull request, the developer has indeed got back to us quick. Looking forward to fix the addressed comments"
$ws.Range("F42").Value = "We were asked to write test cases for the change we suggested. We are happy that we got a response back and are confident that we can integrate our knowledge from Testing course to write beautiful cases covering the improvement. I liked how detailed the developer had replied, it is easy to understand what he is looking for."
$ws.Range("G42").Value = "Feeling good!"
